# ---------------------------------------------------------------
# France National 2023-2024 — corrections + new fixture (row 119)
# Script generated to reproduce the committed diff exactly.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mismatched match data (odds/results had been attached to the
#     wrong fixture row) -------------------------------------------

# Row 104 <= old Row 105 data (F:V)
$ws.Cells.Item(104,6).Value2 = "Sochaux"
$ws.Cells.Item(104,7).Value2 = 2
$ws.Cells.Item(104,8).Value2 = "Cholet"
$ws.Cells.Item(104,9).Value2 = 0
$ws.Cells.Item(104,10).Value2 = 1.53
$ws.Cells.Item(104,11).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(104,12).Value2 = 1.61
$ws.Cells.Item(104,13).Value2 = "10/11/2023 19:17"
$ws.Cells.Item(104,14).Value2 = 3.97
$ws.Cells.Item(104,15).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(104,16).Value2 = 3.92
$ws.Cells.Item(104,17).Value2 = "10/11/2023 19:17"
$ws.Cells.Item(104,18).Value2 = 5.68
$ws.Cells.Item(104,19).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(104,20).Value2 = 5.75
$ws.Cells.Item(104,21).Value2 = "10/11/2023 19:29"
$ws.Cells.Item(104,22).Value2 = "https://www.betexplorer.com/football/france/national/sochaux-cholet/WvqtgGYk/"

# Row 105 <= old Row 104 data (F:V)
$ws.Cells.Item(105,6).Value2 = "Versailles"
$ws.Cells.Item(105,7).Value2 = 6
$ws.Cells.Item(105,8).Value2 = "Nimes"
$ws.Cells.Item(105,9).Value2 = 0
$ws.Cells.Item(105,10).Value2 = 2.16
$ws.Cells.Item(105,11).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(105,12).Value2 = 2.35
$ws.Cells.Item(105,13).Value2 = "10/11/2023 19:21"
$ws.Cells.Item(105,14).Value2 = 3.09
$ws.Cells.Item(105,15).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(105,16).Value2 = 2.87
$ws.Cells.Item(105,17).Value2 = "10/11/2023 19:21"
$ws.Cells.Item(105,18).Value2 = 3.72
$ws.Cells.Item(105,19).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(105,20).Value2 = 3.69
$ws.Cells.Item(105,21).Value2 = "10/11/2023 19:21"
$ws.Cells.Item(105,22).Value2 = "https://www.betexplorer.com/football/france/national/versailles-nimes/0WzGAwbg/"

# Row 107 <= old Row 109 data (F:V)
$ws.Cells.Item(107,6).Value2 = "GOAL FC"
$ws.Cells.Item(107,7).Value2 = 2
$ws.Cells.Item(107,8).Value2 = "Niort"
$ws.Cells.Item(107,9).Value2 = 3
$ws.Cells.Item(107,10).Value2 = 2.45
$ws.Cells.Item(107,11).Value2 = "06/11/2023 18:42"
$ws.Cells.Item(107,12).Value2 = 2.27
$ws.Cells.Item(107,13).Value2 = "10/11/2023 19:29"
$ws.Cells.Item(107,14).Value2 = 3.11
$ws.Cells.Item(107,15).Value2 = "06/11/2023 18:42"
$ws.Cells.Item(107,16).Value2 = 2.93
$ws.Cells.Item(107,17).Value2 = "10/11/2023 19:26"
$ws.Cells.Item(107,18).Value2 = 3.08
$ws.Cells.Item(107,19).Value2 = "06/11/2023 18:42"
$ws.Cells.Item(107,20).Value2 = 3.81
$ws.Cells.Item(107,21).Value2 = "10/11/2023 19:29"
$ws.Cells.Item(107,22).Value2 = "https://www.betexplorer.com/football/france/national/goal-fc-niort/C6vK9cDa/"

# Row 109 <= old Row 107 data (F:V)
$ws.Cells.Item(109,6).Value2 = "Epinal"
$ws.Cells.Item(109,7).Value2 = 4
$ws.Cells.Item(109,8).Value2 = "Avranches"
$ws.Cells.Item(109,9).Value2 = 0
$ws.Cells.Item(109,10).Value2 = 2.73
$ws.Cells.Item(109,11).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(109,12).Value2 = 2.53
$ws.Cells.Item(109,13).Value2 = "10/11/2023 19:29"
$ws.Cells.Item(109,14).Value2 = 3.21
$ws.Cells.Item(109,15).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(109,16).Value2 = 3.22
$ws.Cells.Item(109,17).Value2 = "10/11/2023 19:27"
$ws.Cells.Item(109,18).Value2 = 2.67
$ws.Cells.Item(109,19).Value2 = "04/11/2023 04:43"
$ws.Cells.Item(109,20).Value2 = 2.96
$ws.Cells.Item(109,21).Value2 = "10/11/2023 19:29"
$ws.Cells.Item(109,22).Value2 = "https://www.betexplorer.com/football/france/national/epinal-avranches/2eXW6ecI/"

# Row 112 <= old Row 113 data (F:V)
$ws.Cells.Item(112,6).Value2 = "Chateauroux"
$ws.Cells.Item(112,7).Value2 = 1
$ws.Cells.Item(112,8).Value2 = "Avranches"
$ws.Cells.Item(112,9).Value2 = 2
$ws.Cells.Item(112,10).Value2 = 1.97
$ws.Cells.Item(112,11).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(112,12).Value2 = 1.85
$ws.Cells.Item(112,13).Value2 = "24/11/2023 19:29"
$ws.Cells.Item(112,14).Value2 = 3.38
$ws.Cells.Item(112,15).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(112,16).Value2 = 3.37
$ws.Cells.Item(112,17).Value2 = "24/11/2023 19:29"
$ws.Cells.Item(112,18).Value2 = 3.97
$ws.Cells.Item(112,19).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(112,20).Value2 = 4.73
$ws.Cells.Item(112,21).Value2 = "24/11/2023 19:29"
$ws.Cells.Item(112,22).Value2 = "https://www.betexplorer.com/football/france/national/chateauroux-avranches/KnWkNGsa/"

# Row 113 <= old Row 114 data (F:V)
$ws.Cells.Item(113,6).Value2 = "Cholet"
$ws.Cells.Item(113,7).Value2 = 0
$ws.Cells.Item(113,8).Value2 = "Martigues"
$ws.Cells.Item(113,9).Value2 = 2
$ws.Cells.Item(113,10).Value2 = 3.51
$ws.Cells.Item(113,11).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(113,12).Value2 = 4.34
$ws.Cells.Item(113,13).Value2 = "24/11/2023 19:27"
$ws.Cells.Item(113,14).Value2 = 3.16
$ws.Cells.Item(113,15).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(113,16).Value2 = 3.15
$ws.Cells.Item(113,17).Value2 = "24/11/2023 19:27"
$ws.Cells.Item(113,18).Value2 = 2.21
$ws.Cells.Item(113,19).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(113,20).Value2 = 2.01
$ws.Cells.Item(113,21).Value2 = "24/11/2023 19:27"
$ws.Cells.Item(113,22).Value2 = "https://www.betexplorer.com/football/france/national/cholet-martigues/63Y1KERH/"

# Row 114 <= old Row 115 data (F:V)
$ws.Cells.Item(114,6).Value2 = "Dijon"
$ws.Cells.Item(114,7).Value2 = 3
$ws.Cells.Item(114,8).Value2 = "Red Star"
$ws.Cells.Item(114,9).Value2 = 1
$ws.Cells.Item(114,10).Value2 = 2.2
$ws.Cells.Item(114,11).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(114,12).Value2 = 3.02
$ws.Cells.Item(114,13).Value2 = "24/11/2023 19:26"
$ws.Cells.Item(114,14).Value2 = 3.28
$ws.Cells.Item(114,15).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(114,16).Value2 = 3.18
$ws.Cells.Item(114,17).Value2 = "24/11/2023 19:26"
$ws.Cells.Item(114,18).Value2 = 3.21
$ws.Cells.Item(114,19).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(114,20).Value2 = 2.51
$ws.Cells.Item(114,21).Value2 = "24/11/2023 19:26"
$ws.Cells.Item(114,22).Value2 = "https://www.betexplorer.com/football/france/national/dijon-red-star/KC1oQIct/"

# Row 115 <= old Row 116 data (F:V)
$ws.Cells.Item(115,6).Value2 = "Nancy"
$ws.Cells.Item(115,7).Value2 = 1
$ws.Cells.Item(115,8).Value2 = "Rouen"
$ws.Cells.Item(115,9).Value2 = 0
$ws.Cells.Item(115,10).Value2 = 3.01
$ws.Cells.Item(115,11).Value2 = "20/11/2023 05:12"
$ws.Cells.Item(115,12).Value2 = 3.86
$ws.Cells.Item(115,13).Value2 = "24/11/2023 19:21"
$ws.Cells.Item(115,14).Value2 = 2.84
$ws.Cells.Item(115,15).Value2 = "20/11/2023 05:12"
$ws.Cells.Item(115,16).Value2 = 2.88
$ws.Cells.Item(115,17).Value2 = "24/11/2023 19:21"
$ws.Cells.Item(115,18).Value2 = 2.59
$ws.Cells.Item(115,19).Value2 = "20/11/2023 05:12"
$ws.Cells.Item(115,20).Value2 = 2.28
$ws.Cells.Item(115,21).Value2 = "24/11/2023 19:21"
$ws.Cells.Item(115,22).Value2 = "https://www.betexplorer.com/football/france/national/nancy-rouen/riAJUbkP/"

# Row 116 <= old Row 117 data (F:V)
$ws.Cells.Item(116,6).Value2 = "Nimes"
$ws.Cells.Item(116,7).Value2 = 2
$ws.Cells.Item(116,8).Value2 = "Orleans"
$ws.Cells.Item(116,9).Value2 = 3
$ws.Cells.Item(116,10).Value2 = 2.32
$ws.Cells.Item(116,11).Value2 = "17/11/2023 18:42"
$ws.Cells.Item(116,12).Value2 = 3.16
$ws.Cells.Item(116,13).Value2 = "24/11/2023 19:27"
$ws.Cells.Item(116,14).Value2 = 3.15
$ws.Cells.Item(116,15).Value2 = "17/11/2023 18:42"
$ws.Cells.Item(116,16).Value2 = 2.74
$ws.Cells.Item(116,17).Value2 = "24/11/2023 19:15"
$ws.Cells.Item(116,18).Value2 = 3.27
$ws.Cells.Item(116,19).Value2 = "17/11/2023 18:42"
$ws.Cells.Item(116,20).Value2 = 2.76
$ws.Cells.Item(116,21).Value2 = "24/11/2023 19:27"
$ws.Cells.Item(116,22).Value2 = "https://www.betexplorer.com/football/france/national/nimes-orleans/h0UcLfCB/"

# Row 117 <= old Row 112 data (F:V)
$ws.Cells.Item(117,6).Value2 = "Villefranche"
$ws.Cells.Item(117,7).Value2 = 2
$ws.Cells.Item(117,8).Value2 = "Le Mans"
$ws.Cells.Item(117,9).Value2 = 0
$ws.Cells.Item(117,10).Value2 = 2.33
$ws.Cells.Item(117,11).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(117,12).Value2 = 2.61
$ws.Cells.Item(117,13).Value2 = "24/11/2023 19:24"
$ws.Cells.Item(117,14).Value2 = 3.23
$ws.Cells.Item(117,15).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(117,16).Value2 = 3.04
$ws.Cells.Item(117,17).Value2 = "24/11/2023 19:04"
$ws.Cells.Item(117,18).Value2 = 3.17
$ws.Cells.Item(117,19).Value2 = "19/11/2023 23:12"
$ws.Cells.Item(117,20).Value2 = 3.01
$ws.Cells.Item(117,21).Value2 = "24/11/2023 19:24"
$ws.Cells.Item(117,22).Value2 = "https://www.betexplorer.com/football/france/national/villefranche-le-mans/hQagOdSh/"

# --- Append the new fixture (Niort 4-2 Sochaux) as row 119 ---------
# Copy formatting (styles/number formats) from the last existing row,
# then overwrite with the new match data.
$ws.Range("A118:V118").Copy($ws.Range("A119:V119"))

$ws.Cells.Item(119,1).Value2 = 118
$ws.Cells.Item(119,2).Value2 = "france"
$ws.Cells.Item(119,3).Value2 = "national"
$ws.Cells.Item(119,4).Value2 = "2023-2024"
$ws.Cells.Item(119,5).Value2 = 45257.77083333334
$ws.Cells.Item(119,6).Value2 = "Niort"
$ws.Cells.Item(119,7).Value2 = 4
$ws.Cells.Item(119,8).Value2 = "Sochaux"
$ws.Cells.Item(119,9).Value2 = 2
$ws.Cells.Item(119,10).Value2 = 2.52
$ws.Cells.Item(119,11).Value2 = "20/11/2023 04:42"
$ws.Cells.Item(119,12).Value2 = 2.88
$ws.Cells.Item(119,13).Value2 = "27/11/2023 18:29"
$ws.Cells.Item(119,14).Value2 = 3.3
$ws.Cells.Item(119,15).Value2 = "20/11/2023 04:42"
$ws.Cells.Item(119,16).Value2 = 3.34
$ws.Cells.Item(119,17).Value2 = "27/11/2023 18:29"
$ws.Cells.Item(119,18).Value2 = 2.7
$ws.Cells.Item(119,19).Value2 = "20/11/2023 04:42"
$ws.Cells.Item(119,20).Value2 = 2.52
$ws.Cells.Item(119,21).Value2 = "27/11/2023 18:29"
$ws.Cells.Item(119,22).Value2 = "https://www.betexplorer.com/football/france/national/niort-sochaux/AVsphzld/"
